$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: delete the very first paragraph entirely
#   "SCOALA ................................................"
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(1).Range.Delete()

# ---------------------------------------------------------------------------
# Step 2: rewrite the "AVIZAT," paragraph (now paragraph 1): the run holding
# 24 leading spaces is dropped and replaced by two <w:tab/> placed right
# before the "AVIZAT, " run.
# ---------------------------------------------------------------------------
$pAvizat = $d.Paragraphs.Item(1).Range
$pAvizat.Collapse(1)
$xmlAvizat = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:before="0" w:after="0"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ro-RO"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:cs="Times New Roman" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ro-RO"/></w:rPr><w:tab/><w:tab/><w:t xml:space="preserve">AVIZAT, </w:t></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pAvizat.InsertXML($xmlAvizat)

# ---------------------------------------------------------------------------
# Step 3: rewrite the "DIRECTOR" paragraph (now paragraph 2): split the
# "DIRECTOR" + 13 tabs run into "DIRECTO" / "R " / "SCOALA " /
# "PROFESIOANALA SPECIALA ..." runs, dropping every tab.
# The leading 22-space run stays untouched.
# ---------------------------------------------------------------------------
$pDirector = $d.Paragraphs.Item(2).Range
$pDirector.Collapse(1)
$xmlDirector = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:before="0" w:after="0"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ro-RO"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:cs="Times New Roman" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve">                      </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:cs="Times New Roman" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ro-RO"/></w:rPr><w:t>DIRECTO</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:cs="Times New Roman" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve">R </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:cs="Times New Roman" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve">ȘCOALA </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:cs="Times New Roman" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ro-RO"/></w:rPr><w:t>PROFESIOANALĂ SPECIALĂ „ION TEODORESCU”</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pDirector.InsertXML($xmlDirector)

# ---------------------------------------------------------------------------
# Step 4: the following (previously empty) paragraph becomes "Profesor,"
# preceded by two tabs, and two brand-new paragraphs are inserted right
# after it (a 24-space line, then a blank formatted line) before the
# Algerian-font tab paragraph that must stay untouched.
# ---------------------------------------------------------------------------
$pProfesor = $d.Paragraphs.Item(3).Range
$pProfesor.Collapse(1)
$xmlProfesor = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:before="0" w:after="0"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ro-RO"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:cs="Times New Roman" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ro-RO"/></w:rPr><w:tab/><w:tab/><w:t>Profesor,</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:before="0" w:after="0"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ro-RO"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:cs="Times New Roman" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve">                        </w:t></w:r>' + `
  '</w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:before="0" w:after="0"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ro-RO"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:cs="Times New Roman" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="ro-RO"/></w:rPr></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pProfesor.InsertXML($xmlProfesor)

Write-Output "done"
